$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 73

# Column A holds a date-formatted label ("01-07-2021") stored as plain text
# (shared string) in the source file, matching the other "Serie" cells
# above it. Assigning the literal text directly would make Excel's
# automatic date recognition convert it into a date serial + date number
# format. Routing it through a formula (so the smart-type parser is never
# invoked) and then converting that formula to its literal value via
# copy / paste-special-values keeps it as plain text without touching the
# cell's style.
$ws.Cells.Item($row, 1).Formula = '="01-07-2021"'
$ws.Cells.Item($row, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 2).Value = 60.7
$ws.Cells.Item($row, 3).Value = 49.7
$ws.Cells.Item($row, 4).Value = 52.1
$ws.Cells.Item($row, 5).Value = 94.2
$ws.Cells.Item($row, 6).Value = 92
$ws.Cells.Item($row, 7).Value = 82.8
$ws.Cells.Item($row, 8).Value = 75.2
$ws.Cells.Item($row, 9).Value = 50
$ws.Cells.Item($row, 10).Value = 55.9
$ws.Cells.Item($row, 11).Value = 77.5
$ws.Cells.Item($row, 12).Value = 91.5
$ws.Cells.Item($row, 13).Value = 87.40000000000001
$ws.Cells.Item($row, 14).Value = 93.8
$ws.Cells.Item($row, 15).Value = 61.7
